# se agrego titulo a la presentacion
#
# Slide 1 has two empty placeholders (a center-title and a subtitle).
# Fill them in with the title / subtitle text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1 = "Título 1" (ctrTitle placeholder)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Gracias a la compañera "

# Shape 2 = "Subtítulo 2" (subTitle placeholder)
$subtitle = $s.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = "Por fin "
